$d = $word.ActiveDocument

# 1. Merge "(i) state the identity..." split across proofErr spell-check runs
$d.Content.Find.Execute(
    "(i) state the identity of the client and any intended users, by name or type;",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "(i) state the identity of the client and any intended users, by name or type;", 2)

# 2. "Fee Simple" -> "Leased Fee" (Purpose of the Appraisal section)
$d.Content.Find.Execute(
    "Fee Simple", $true, $false, $false, $false, $false, $true, 1, $false,
    "Leased Fee", 2)

# 3. Merge "The most probable price..." definition split by gramStart/gramEnd
$d.Content.Find.Execute(
    "The most probable price which a property should bring in a competitive and open market.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "The most probable price which a property should bring in a competitive and open market.", 2)

$d.Content.Find.Execute(
    "Conditions requisite to a fair sale,",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Conditions requisite to a fair sale,", 2)

# 4. Merge "Exposure time as used in this appraisal report is defined as:"
$d.Content.Find.Execute(
    "Exposure time as used in this appraisal report is defined as:",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Exposure time as used in this appraisal report is defined as:", 2)

# 5 & 7. Merge "Based upon the market data from the sales of comparable properties in the market my conclusion of exposure time follows:" (occurs twice)
$d.Content.Find.Execute(
    "Based upon the market data from the sales of comparable properties in the market my conclusion of exposure time follows:",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Based upon the market data from the sales of comparable properties in the market my conclusion of exposure time follows:", 2)

# 6. Merge "subsequent to" split phrase
$d.Content.Find.Execute(
    "The time it takes an interest in real property to sell on the market subsequent to the date of appraisal",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "The time it takes an interest in real property to sell on the market subsequent to the date of appraisal", 2)
